# Regenerate column G ("K" = strikeouts) values for the save_data sheet.
# Previously column G held the "Strike#" (count of strike pitches); it has
# been regenerated to hold the actual strikeout count (K) per the commit
# message: "regen save_data to use K instead of Strike#".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 6
    3  = 10
    4  = 7
    5  = 8
    6  = 7
    7  = 5
    8  = 9
    9  = 2
    10 = 8
    11 = 5
    12 = 7
    13 = 11
    14 = 11
    15 = 12
    16 = 10
    17 = 11
    18 = 13
    19 = 6
    20 = 6
    21 = 2
    22 = 10
    23 = 7
    24 = 5
    25 = 11
    26 = 4
    27 = 12
    28 = 5
    29 = 9
    30 = 4
    31 = 5
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
